# Auto commit at 2026-01-15  8:13:55.97
# Updates the Metrics sheet values (which cascade via formulas into the
# "today" sheet), and restores the two sheet selections that Excel
# persists in sheetView/selection when the workbook was last saved.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refresh the metric values ------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value  = 266183.55
$wsMetrics.Range("B3").Value  = 197691.99000000002
$wsMetrics.Range("B4").Value  = 68693.110000000015
$wsMetrics.Range("B5").Value  = 10848
$wsMetrics.Range("B6").Value  = 5902054.2799999984
$wsMetrics.Range("B7").Value  = 4968409.62
$wsMetrics.Range("B8").Value  = 1732784.93
$wsMetrics.Range("B9").Value  = 231125
$wsMetrics.Range("B10").Value = 34367435.269999996
$wsMetrics.Range("B11").Value = 32243684.780000001
$wsMetrics.Range("B12").Value = 12014506.970000001
$wsMetrics.Range("B13").Value = 1328755

# Restore the cursor/selection saved in the Metrics sheet view.
$wsMetrics.Activate()
$wsMetrics.Range("D14").Select()

# --- today sheet: move the saved selection ---------------------------
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate()
$wsToday.Range("E6").Select()
